# [ADD] - Ajout de l'introduction dans le rapport
$d = $word.ActiveDocument

# 1) Insert a new empty paragraph at the very start of the body
#    (jc=both, rFonts Aller) - matches the existing first paragraph's
#    justification and uses the "Aller" heading font family.
$startRng = $d.Range(0, 0)
$startRng.InsertParagraphBefore()
$firstP = $d.Paragraphs(1)
$firstP.Alignment = 3
$firstP.Range.Font.Name = "Aller"

# 2) Locate the "Introduction" title paragraph (Titre_Partie style) and:
#    - drop the stale <w:lastRenderedPageBreak/>
#    - split "Introduction" into "Introductio" + "n" around the bookmark
$introIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd() -eq "Introduction" -and $p.Style.NameLocal -eq "Titre_Partie") {
        $introIdx = $i
        break
    }
}

$introP = $d.Paragraphs($introIdx)
$introRng = $introP.Range
$introRng.Find.Execute("Introduction", $true, $true, $false, $false, $false, $true, 1, $false, "Introduction", 2) | Out-Null

$bm = $d.Bookmarks("_Toc128474534")
$splitRange = $d.Range($bm.Start, $bm.Start + 11)
$d.Bookmarks.Add("_Toc128474534", $splitRange)

# 3) Insert the new introductory paragraphs between the "Introduction"
#    title and the next title ("Editeur de monde").
$targetIdx = $introIdx + 1

$d.Range($d.Paragraphs($targetIdx).Range.Start, $d.Paragraphs($targetIdx).Range.Start).InsertParagraphBefore()
$para1 = $d.Paragraphs($targetIdx)
$para1.Style = "Paragraphe"
$pos = $para1.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Nous avons entrepris la conception d’un projet ambitieux appelé Starlyze, qui consiste à développer un jeu de plateforme multijoueur en ligne. L’objectif est de créer un monde virtuel dans lequel chaque joueur peut contrôler un personnage et interagir avec d’autres joueurs de la même partie. Le jeu est conçu pour être immersif, passionnant et stimulant pour les joueurs. ")
$targetIdx = $targetIdx + 1

$d.Range($d.Paragraphs($targetIdx).Range.Start, $d.Paragraphs($targetIdx).Range.Start).InsertParagraphBefore()
$para2 = $d.Paragraphs($targetIdx)
$para2.Style = "Paragraphe"
$pos = $para2.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Le concept de base du jeu est simple : chaque joueur peut créer son propre monde et y jouer avec d’autres joueurs. Cela signifie qu’il n’y aura pas de monde prédéfinies, et chaque joueur aura la liberté de créer son propre monde avec ses propres niveaux, décors, obstacles et ennemis. Ainsi, chaque partie de Starlyze sera unique et offrira une expérience de jeu différente à chaque fois.")
$targetIdx = $targetIdx + 1

$d.Range($d.Paragraphs($targetIdx).Range.Start, $d.Paragraphs($targetIdx).Range.Start).InsertParagraphBefore()
$para3 = $d.Paragraphs($targetIdx)
$para3.Style = "Paragraphe"
$pos = $para3.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Le projet est divisé en trois parties principales : l’éditeur, le serveur et le client. L’éditeur permet aux joueur de créer leur propre monde et de le personnaliser selon leurs préférences. Les joueurs peuvent ajouter des objets, des niveaux et des décors à leur monde. Le serveur est la partie du jeu qui gère la mise en relation des joueurs et la gestion de la partie. Il est responsable de la création de parties, de la gestion des joueurs et de la communication avec les joueurs. Le serveur est conçu afin d’être performant et fiable, afin d’offrir une expérience de jeu fluide et sans interruption. Enfin, le client est l’application qui permet aux joueurs de visualiser et de jouer au jeu. Il est conçu pour être facile à utiliser et à comprendre, avec une interface utilisateur intuitive et ergonomique. Le client affiche la partie en temps réel et permet aux joueurs de contrôler leur personnage.")
$targetIdx = $targetIdx + 1

$d.Range($d.Paragraphs($targetIdx).Range.Start, $d.Paragraphs($targetIdx).Range.Start).InsertParagraphBefore()
$para4 = $d.Paragraphs($targetIdx)
$para4.Style = "Paragraphe"
$para4.Format.FirstLineIndent = 0
$pos = $para4.Range.End - 1
$d.Range($pos, $pos).InsertAfter("De plus, les applications « serveur » et « client » communiquent tout au long de la partie via différents protocoles afin de garantir la fluidité de l’expérience de jeu et une interaction en temps réel entre les joueurs.")
$targetIdx = $targetIdx + 1

$d.Range($d.Paragraphs($targetIdx).Range.Start, $d.Paragraphs($targetIdx).Range.Start).InsertParagraphBefore()
$para5 = $d.Paragraphs($targetIdx)
$para5.Style = "Paragraphe"
$pos = $para5.Range.End - 1
$d.Range($pos, $pos).InsertAfter("En somme, Starlyze est un projet conséquent qui vise à offrir une expérience de jeu unique et immersive à tous les joueurs. Avec sa conception intuitive et son gameplay stimulant, ce jeu promet d’être une référence dans le monde des jeux de plateforme multijoueur en ligne.")

# 4) Remove the other stale <w:lastRenderedPageBreak/> further down, on the
#    "On rappelle que pour un monde ..." paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("On rappelle que pour un monde")) {
        $rng = $p.Range
        $fullText = $rng.Text.TrimEnd()
        $rng.Find.Execute($fullText, $true, $true, $false, $false, $false, $true, 1, $false, $fullText, 2) | Out-Null
        break
    }
}
